$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exceptions")

$ws.Range("A8").Value = "0x0006"
$ws.Range("B8").Value = "DataBase.cs"
$ws.Range("C8").Value = "Unable to close the DB connection"

$ws.Range("C8").Select()
